# "read only defined columns"
# Populate row 3 (A3:E3) the same way row 2 (A2:E2) is populated, and
# move the active selection to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "The same"
$ws.Range("B3").Value = 4534567
$ws.Range("C3").Value = " a table"
$ws.Range("D3").Value = 0
$ws.Range("E3").Formula = '=IF(AND($A3<>"",$B3<>"",$D3<>""),"Filled","Not filled")'

$ws.Range("B3").Select() | Out-Null
